$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so values like
# "584.30" or "0.520" keep their exact formatting (no numeric coercion).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.666.50"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.282.95"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.30"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.17"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.419"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.856.08"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.38"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.708.42"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.295.01"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.88"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.60"
$ws.Range("E20").Value = "  +4.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.75"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.75"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.520"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("E26").Value = "  +4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.78"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.76"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.10"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.16"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.38"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.01"
$ws.Range("E37").Value = "  +8.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.826"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.62"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.54"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.61"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.43"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0693"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.651.23"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.47"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0283"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.06"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.37"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -0.82%  "
